# Apply weekly data update: insert a new data row at row 194 (shifting the
# existing rows 194-291 down to 195-292), and populate the new row with the
# latest week's values while keeping the other (unchanged) columns equal to
# the values that used to be in the old row 194.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the values from the current row 194 that remain identical in the
# new row 194 (everything except Fecha, Volumen, Precio minimo, Precio
# maximo, Precio promedio ponderado, Origen and Precio $/Kg).
$mercadoId  = $ws.Range("A194").Value2
$mercado    = $ws.Range("B194").Value2
$region     = $ws.Range("C194").Value2
$codreg     = $ws.Range("E194").Value2
$categoriaId = $ws.Range("F194").Value2
$categoria  = $ws.Range("G194").Value2
$variedad   = $ws.Range("H194").Value2
$calidad    = $ws.Range("I194").Value2
$unidad     = $ws.Range("N194").Value2
$kgOUnid    = $ws.Range("Q194").Value2
$clasif     = $ws.Range("R194").Value2

# Keep a reference to the current date cell's number format so the new row
# can reuse the same (date) formatting.
$dateFormat = $ws.Range("D194").NumberFormat

# Insert a new blank row at 194; this pushes the old rows 194-291 down to
# 195-292 and extends the used range to row 292.
$ws.Rows("194:194").Insert()

# Re-write the columns that stay the same as before.
$ws.Range("A194").Value = $mercadoId
$ws.Range("B194").Value = $mercado
$ws.Range("C194").Value = $region
$ws.Range("E194").Value = $codreg
$ws.Range("F194").Value = $categoriaId
$ws.Range("G194").Value = $categoria
$ws.Range("H194").Value = $variedad
$ws.Range("I194").Value = $calidad
$ws.Range("N194").Value = $unidad
$ws.Range("Q194").Value = $kgOUnid
$ws.Range("R194").Value = $clasif

# New values for the newly inserted row.
$ws.Range("D194").Value = 44873
$ws.Range("D194").NumberFormat = $dateFormat
$ws.Range("J194").Value = 270
$ws.Range("K194").Value = 15000
$ws.Range("L194").Value = 16000
$ws.Range("M194").Value = 15556
$ws.Range("O194").Value = "Región Metropolitana"
$ws.Range("P194").Value = 778
